$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "29÷2=14, 1"  # was: 94÷5=18, 4
$t.Cell(1, 2).Range.Text = "64÷9=7, 1"  # was: 75÷8=9, 3
$t.Cell(1, 3).Range.Text = "42÷5=8, 2"  # was: 44÷4=11, 0
$t.Cell(1, 4).Range.Text = "12÷3=4, 0"  # was: 99÷8=12, 3
$t.Cell(1, 5).Range.Text = "42÷5=8, 2"  # was: 71÷9=7, 8

$t.Cell(5, 1).Range.Text = "73÷6=12, 1"  # was: 15÷6=2, 3
$t.Cell(5, 2).Range.Text = "89÷5=17, 4"  # was: 54÷2=27, 0
$t.Cell(5, 3).Range.Text = "67÷3=22, 1"  # was: 34÷5=6, 4
$t.Cell(5, 4).Range.Text = "39÷5=7, 4"  # was: 19÷4=4, 3
$t.Cell(5, 5).Range.Text = "70÷4=17, 2"  # was: 39÷9=4, 3

$t.Cell(9, 1).Range.Text = "82÷7=11, 5"  # was: 52÷4=13, 0
$t.Cell(9, 2).Range.Text = "69÷8=8, 5"  # was: 33÷3=11, 0
$t.Cell(9, 3).Range.Text = "12÷4=3, 0"  # was: 41÷9=4, 5
$t.Cell(9, 4).Range.Text = "34÷7=4, 6"  # was: 45÷7=6, 3
$t.Cell(9, 5).Range.Text = "43÷6=7, 1"  # was: 59÷8=7, 3

$t.Cell(13, 1).Range.Text = "31÷7=4, 3"  # was: 63÷2=31, 1
$t.Cell(13, 2).Range.Text = "24÷6=4, 0"  # was: 61÷2=30, 1
$t.Cell(13, 3).Range.Text = "92÷3=30, 2"  # was: 61÷8=7, 5
$t.Cell(13, 4).Range.Text = "70÷7=10, 0"  # was: 29÷8=3, 5
$t.Cell(13, 5).Range.Text = "30÷2=15, 0"  # was: 17÷5=3, 2

$t.Cell(17, 1).Range.Text = "67÷7=9, 4"  # was: 84÷9=9, 3
$t.Cell(17, 2).Range.Text = "21÷2=10, 1"  # was: 29÷2=14, 1
$t.Cell(17, 3).Range.Text = "81÷6=13, 3"  # was: 46÷8=5, 6
$t.Cell(17, 4).Range.Text = "39÷6=6, 3"  # was: 58÷5=11, 3
$t.Cell(17, 5).Range.Text = "59÷7=8, 3"  # was: 69÷2=34, 1
